$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new timelog entry: week 3, 2nd class
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A9").Value = 45919
$ws.Range("B9").Value = "class"
$ws.Range("C9").Value = "9:30am"
$ws.Range("D9").Value = "12:25pm"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 55

# Grow Table1 so it covers the newly added row
$ws.ListObjects("Table1").Resize($ws.Range("A1:F9"))

# Update selection to mirror the saved workbook state
$ws.Range("G10").Select()
